$wb = $excel.ActiveWorkbook

# addListItem sheet: update the eli_text (A2) and eli_code_type (D2) values
$wsAdd = $wb.Worksheets.Item("addListItem")
$wsAdd.Range("A2").Value = "UserELIp"
$wsAdd.Range("D2").Value = "ADLILC.8847"

# createUser sheet: bump the numeric id in A2; dependent formulas in B2/F2 recalc automatically
$wsUser = $wb.Worksheets.Item("createUser")
$wsUser.Range("A2").Value = 150

$wb.Save()
